# Auto-generated Excel COM-interop script to apply scheduled Sheets update
# Updates currentAveragePrice / price / profit columns (H-N) across multiple
# crafting-leve worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 40
$ws.Range("H40").Value = 7824.857
$ws.Range("I40").Value = 6194
$ws.Range("K40").Value = 6194
$ws.Range("M40").Value = -6019
# row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# row 116
$ws.Range("H116").Value = 4714.364
$ws.Range("I116").Value = 4715.9
$ws.Range("J116").Value = 4699
$ws.Range("K116").Value = 4715.9
$ws.Range("L116").Value = 4699
$ws.Range("M116").Value = -1273.9
$ws.Range("N116").Value = -11583
# row 125
$ws.Range("H125").Value = 934
$ws.Range("I125").Value = 932
$ws.Range("J125").Value = 936
$ws.Range("K125").Value = 8388
$ws.Range("L125").Value = 8424
$ws.Range("M125").Value = -5928
$ws.Range("N125").Value = -13344
# row 141
$ws.Range("H141").Value = 1049.5454
$ws.Range("I141").Value = 1049.5454
$ws.Range("K141").Value = 3148.6362
$ws.Range("M141").Value = 2031.3638

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 1124.9
$ws.Range("I2").Value = 1124.9
$ws.Range("K2").Value = 1124.9
$ws.Range("M2").Value = -1011.9
# row 32
$ws.Range("H32").Value = 949.6070999999999
$ws.Range("I32").Value = 868.8077
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 868.8077
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -581.8077
$ws.Range("N32").Value = -2574
# row 45
$ws.Range("H45").Value = 3018.4375
$ws.Range("I45").Value = 2018.8889
$ws.Range("K45").Value = 2018.8889
$ws.Range("M45").Value = -1641.8889
# row 63
$ws.Range("H63").Value = 5443.5713
$ws.Range("I63").Value = 2866.6667
$ws.Range("J63").Value = 7376.25
$ws.Range("K63").Value = 2866.6667
$ws.Range("L63").Value = 7376.25
$ws.Range("M63").Value = -2180.6667
$ws.Range("N63").Value = -8748.25
# row 66
$ws.Range("H66").Value = 5443.5713
$ws.Range("I66").Value = 2866.6667
$ws.Range("J66").Value = 7376.25
$ws.Range("K66").Value = 14333.3335
$ws.Range("L66").Value = 36881.25
$ws.Range("M66").Value = -10901.3335
$ws.Range("N66").Value = -43745.25
# row 97
$ws.Range("H97").Value = 1061
$ws.Range("I97").Value = 1037.0769
$ws.Range("K97").Value = 1037.0769
$ws.Range("M97").Value = -541.0769
# row 102
$ws.Range("H102").Value = 3200.9285
$ws.Range("I102").Value = 859
$ws.Range("J102").Value = 5542.857
$ws.Range("K102").Value = 859
$ws.Range("L102").Value = 5542.857
$ws.Range("M102").Value = 763
$ws.Range("N102").Value = -8786.857
# row 116
$ws.Range("H116").Value = 1124.9
$ws.Range("I116").Value = 1124.9
$ws.Range("K116").Value = 1124.9
$ws.Range("M116").Value = 1169.1
# row 132
$ws.Range("H132").Value = 1341.5714
$ws.Range("I132").Value = 1602.8
$ws.Range("J132").Value = 688.5
$ws.Range("K132").Value = 4808.4
$ws.Range("L132").Value = 2065.5
$ws.Range("M132").Value = -2278.4
$ws.Range("N132").Value = -7125.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 1124.9
$ws.Range("I3").Value = 1124.9
$ws.Range("K3").Value = 1124.9
$ws.Range("M3").Value = -1010.9
# row 22
$ws.Range("H22").Value = 270.42856
$ws.Range("I22").Value = 278.2
$ws.Range("J22").Value = 251
$ws.Range("K22").Value = 278.2
$ws.Range("L22").Value = 251
$ws.Range("M22").Value = -105.2
$ws.Range("N22").Value = -597
# row 94
$ws.Range("H94").Value = 853.4
$ws.Range("I94").Value = 630.9231
$ws.Range("K94").Value = 630.9231
$ws.Range("M94").Value = -179.9231
# row 97
$ws.Range("H97").Value = 9233
$ws.Range("I97").Value = 9233
$ws.Range("K97").Value = 9233
$ws.Range("M97").Value = -8242
# row 99
$ws.Range("H99").Value = 1274.6666
$ws.Range("I99").Value = 1274.6666
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1274.6666
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 223.3334
$ws.Range("N99").ClearContents()
# row 107
$ws.Range("H107").Value = 4451
$ws.Range("I107").Value = 4697.1
$ws.Range("J107").Value = 1990
$ws.Range("K107").Value = 4697.1
$ws.Range("L107").Value = 1990
$ws.Range("M107").Value = -2777.1
$ws.Range("N107").Value = -5830

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 105
$ws.Range("H105").Value = 3590.6667
$ws.Range("I105").Value = 3590.6667
$ws.Range("K105").Value = 3590.6667
$ws.Range("M105").Value = -1843.6667
# row 107
$ws.Range("H107").Value = 835.15
$ws.Range("I107").Value = 334.64285
$ws.Range("K107").Value = 334.64285
$ws.Range("M107").Value = 1585.35715
# row 132
$ws.Range("H132").Value = 3161
$ws.Range("I132").Value = 3326
$ws.Range("K132").Value = 9978
$ws.Range("M132").Value = -7448

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 6
$ws.Range("H6").Value = 262.33334
$ws.Range("I6").Value = 74.8
$ws.Range("K6").Value = 224.4
$ws.Range("M6").Value = -111.4
# row 33
$ws.Range("H33").Value = 37.5
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
# row 46
$ws.Range("H46").Value = 5479.6
$ws.Range("I46").Value = 799.3333
$ws.Range("K46").Value = 2397.9999
$ws.Range("M46").Value = -2306.9999
# row 86
$ws.Range("H86").Value = 565.44446
$ws.Range("I86").Value = 561.125
$ws.Range("K86").Value = 1683.375
$ws.Range("M86").Value = -497.375
# row 89
$ws.Range("H89").Value = 565.44446
$ws.Range("I89").Value = 561.125
$ws.Range("K89").Value = 5050.125
$ws.Range("M89").Value = 877.875
# row 92
$ws.Range("H92").Value = 822.8889
$ws.Range("I92").Value = 700.8570999999999
$ws.Range("K92").Value = 2102.5713
$ws.Range("M92").Value = -854.5712999999996
# row 131
$ws.Range("H131").Value = 1554
$ws.Range("J131").Value = 1997.25
$ws.Range("L131").Value = 5991.75
$ws.Range("N131").Value = -16071.75

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 3
$ws.Range("H3").Value = 23676676
$ws.Range("I3").Value = 20864754
$ws.Range("K3").Value = 20864754
$ws.Range("M3").Value = -20864638
# row 70
$ws.Range("H70").Value = 3448.5
$ws.Range("I70").Value = 3448.5
$ws.Range("K70").Value = 3448.5
$ws.Range("M70").Value = -3178.5
# row 73
$ws.Range("H73").Value = 3448.5
$ws.Range("I73").Value = 3448.5
$ws.Range("K73").Value = 3448.5
$ws.Range("M73").Value = -2512.5
# row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
# row 112
$ws.Range("H112").Value = 25000
$ws.Range("J112").Value = 25000
$ws.Range("L112").Value = 25000
$ws.Range("N112").Value = -27216
# row 113
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 1800
$ws.Range("I16").Value = 1800
$ws.Range("K16").Value = 1800
$ws.Range("M16").Value = -1630
# row 22
$ws.Range("H22").Value = 1550.5
$ws.Range("I22").Value = 1246.75
$ws.Range("J22").Value = 1854.25
$ws.Range("K22").Value = 1246.75
$ws.Range("L22").Value = 1854.25
$ws.Range("M22").Value = -951.75
$ws.Range("N22").Value = -2444.25
# row 27
$ws.Range("H27").Value = 1550.5
$ws.Range("I27").Value = 1246.75
$ws.Range("J27").Value = 1854.25
$ws.Range("K27").Value = 1246.75
$ws.Range("L27").Value = 1854.25
$ws.Range("M27").Value = -1139.75
$ws.Range("N27").Value = -2068.25
# row 46
$ws.Range("H46").Value = 6162.722
$ws.Range("I46").Value = 4388.778
$ws.Range("K46").Value = 4388.778
$ws.Range("M46").Value = -4200.778
# row 82
$ws.Range("H82").Value = 3976.0908
$ws.Range("I82").Value = 497.25
$ws.Range("K82").Value = 497.25
$ws.Range("M82").Value = -136.25
# row 85
$ws.Range("H85").Value = 3976.0908
$ws.Range("I85").Value = 497.25
$ws.Range("K85").Value = 497.25
$ws.Range("M85").Value = 750.75
# row 93
$ws.Range("H93").Value = 1358.0714
$ws.Range("I93").Value = 1490.5555
$ws.Range("J93").Value = 1119.6
$ws.Range("K93").Value = 1490.5555
$ws.Range("L93").Value = 1119.6
$ws.Range("M93").Value = -242.5554999999999
$ws.Range("N93").Value = -3615.6
# row 136
$ws.Range("H136").Value = 3147.6667
$ws.Range("I136").Value = 2221.75
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 6665.25
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = -4115.25
$ws.Range("N136").Value = -20098.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 10080.2
$ws.Range("I62").Value = 8000.6665
$ws.Range("J62").Value = 10971.429
$ws.Range("K62").Value = 8000.6665
$ws.Range("L62").Value = 10971.429
$ws.Range("M62").Value = -7376.6665
$ws.Range("N62").Value = -12219.429
# row 65
$ws.Range("H65").Value = 10080.2
$ws.Range("I65").Value = 8000.6665
$ws.Range("J65").Value = 10971.429
$ws.Range("K65").Value = 40003.3325
$ws.Range("L65").Value = 54857.145
$ws.Range("M65").Value = -36883.3325
$ws.Range("N65").Value = -61097.145
# row 100
$ws.Range("H100").Value = 1289.8182
$ws.Range("I100").Value = 1289.8182
$ws.Range("K100").Value = 2579.6364
$ws.Range("M100").Value = -2038.6364
# row 126
$ws.Range("H126").Value = 7986.1665
$ws.Range("I126").Value = 7968.8
$ws.Range("K126").Value = 23906.4
$ws.Range("M126").Value = -21436.4

Write-Host "Applied scheduled Sheets update."